$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# New table contents (old rows 7-15 shifted up to become rows 1-9;
# row 1 keeps its existing header style, rows 10-15 are removed).
$data = @(
    @("CUSTO DESENVOLVIMENTO",     "R$ 68.563,94",    "R$ 0,00",        "R$ 68.563,94",    "R$ 68.563,94",    "100,00 %"),
    @("EMBALAGENS",                "R$ 47.974,08",    "R$ 52.522,04",   "R$ 100.496,12",   "R$ 160.000,00",   "62,81 %"),
    @("DESP. INDUSTRIAL",          "R$ 40.002,66",    "R$ 100.281,50",  "R$ 140.284,16",   "R$ 470.000,00",   "29,85 %"),
    @("SERVICOS DE QUALIDADE",     "R$ 38.449,82",    "R$ 0,00",        "R$ 38.449,82",    "R$ 45.000,00",    "85,44 %"),
    @("CUSTO COM DESENVOLVIMENTO", "R$ 8.301,08",     "R$ 0,00",        "R$ 8.301,08",     "R$ 8.301,08",     "100,00 %"),
    @("FERRAMENTARIA/MAN FR",      "R$ 7.268,32",     "R$ 0,00",        "R$ 7.268,32",     "R$ 35.000,00",    "20,77 %"),
    @("MATERIAL QUALIDADE",        "R$ 110,37",       "R$ 4.033,80",    "R$ 4.144,17",     "R$ 45.000,00",    "9,21 %"),
    @("ENERGIA ELETRICA",          "R$ 0,00",         "R$ 0,00",        "R$ 0,00",         "R$ 519.000,00",   "0 %"),
    @("Total Geral",               "R$ 1.767.242,29", "R$ 534.600,42",  "R$ 2.301.842,71", "R$ 3.835.922,27", "60,01 %")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $cell = $ws.Cells.Item($r, $c)
        $text = $row[$j]
        if ($text -match '^[0-9]+ %$') {
            # Values like "0 %" would otherwise be auto-detected as a
            # percentage number by Excel; force text storage, then drop
            # back to the Normal style so no stray number format sticks.
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.Style = "Normal"
        } else {
            $cell.Value = $text
        }
    }
}

# Remove the now-obsolete trailing rows (old rows 10-15), shrinking the
# used range down to A1:F9.
$ws.Range("A10:F15").Delete()
